# The shape in this test document carries <a:ln> line formatting that is
# not relevant to what the document is actually testing (effect extent
# margins). Strip the <a:ln>...</a:ln> block from the shape's <wps:spPr>
# while leaving everything else (gradient fill, geometry, style, etc.)
# untouched.

$d = $word.ActiveDocument

# Pull the live WordprocessingML for the document so we can splice out
# just the <a:ln> element and feed the result back in verbatim - this
# avoids disturbing any other markup/formatting.
$bodyXml = $d.Content.WordOpenXML

# Locate the shape's <a:ln ...>...</a:ln> (opening tag has attributes,
# e.g. cap="rnd" w="57240", so match on "<a:ln " to avoid also matching
# the unrelated <a:lin ang="..."/> gradient-direction element).
$startTag = "<a:ln "
$endTag = "</a:ln>"

$i1 = $bodyXml.IndexOf($startTag)
if ($i1 -lt 0) {
    throw "a:ln element was not found - nothing to remove"
}
$i2 = $bodyXml.IndexOf($endTag, $i1)
if ($i2 -lt 0) {
    throw "a:ln element was not properly closed"
}
$i2 = $i2 + $endTag.Length

$newXml = $bodyXml.Substring(0, $i1) + $bodyXml.Substring($i2)

$d.Content.InsertXML($newXml)
